$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @(5,7,8,9,10,11,13,14,15,16,17,18,19,20)

$data = @(
    @(3,8.191447666666667,24.574343,0.185794284429433,0.185794284429433,3,8.131233999999999,24.393702,0.02090995573015822,0.02090995573015823,66.60657777642066,599.459199987786,0.003884950262335869,0.003884950262335869),
    @(3,8.191447666666667,24.574343,0.185794284429433,0.185794284429433,3,243.3763986666667,730.1291960000001,0.625857000534647,0.6258570005346471,1993.605032979803,17942.44529681823,0.116280653569486,0.116280653569486),
    @(3,8.191447666666667,24.574343,0.185794284429433,0.185794284429433,3,103.9426383333333,311.827915,0.2672947262403034,0.2672947262403035,851.4406822427605,7662.966140184845,0.04966183239357835,0.04966183239357837),
    @(3,8.191447666666667,24.574343,0.185794284429433,0.185794284429433,3,33.41874933333333,100.256248,0.08593831749489127,0.08593831749489128,273.7479362494515,2463.731426245064,0.01596684820403274,0.01596684820403274),
    @(3,17.317702,51.95310600000001,0.3927913821808575,0.3927913821808576,3,8.131233999999999,24.393702,0.02090995573015822,0.02090995573015823,140.814287304268,1267.328585738412,0.00821325041258939,0.008213250412589393),
    @(3,17.317702,51.95310600000001,0.3927913821808575,0.3927913821808576,3,243.3763986666667,730.1291960000001,0.625857000534647,0.6258570005346471,4214.719945942531,37932.47951348279,0.2458312362875697,0.2458312362875698),
    @(3,17.317702,51.95310600000001,0.3927913821808575,0.3927913821808576,3,103.9426383333333,311.827915,0.2672947262403034,0.2672947262403035,1800.047635750444,16200.42872175399,0.1049910649695827,0.1049910649695827),
    @(3,17.317702,51.95310600000001,0.3927913821808575,0.3927913821808576,3,33.41874933333333,100.256248,0.08593831749489127,0.08593831749489128,578.7359421673654,5208.623479506288,0.03375583051111571,0.03375583051111572),
    @(3,12.35128266666667,37.053848,0.2801455637905346,0.2801455637905346,3,8.131233999999999,24.393702,0.02090995573015822,0.02090995573015823,100.4311695628107,903.880526065296,0.005857831336860294,0.005857831336860295),
    @(3,12.35128266666667,37.053848,0.2801455637905346,0.2801455637905346,3,243.3763986666667,730.1291960000001,0.625857000534647,0.6258570005346471,3006.010694327357,27054.09624894621,0.1753310622670316,0.1753310622670316),
    @(3,12.35128266666667,37.053848,0.2801455637905346,0.2801455637905346,3,103.9426383333333,311.827915,0.2672947262403034,0.2672947262403035,1283.824907174102,11554.42416456692,0.0748814317808264,0.07488143178082642),
    @(3,12.35128266666667,37.053848,0.2801455637905346,0.2801455637905346,3,33.41874933333333,100.256248,0.08593831749489127,0.08593831749489128,412.7644193824783,3714.879774442304,0.02407523840581628,0.02407523840581628),
    @(3,6.228371,18.685113,0.1412687695991749,0.1412687695991749,3,8.131233999999999,24.393702,0.02090995573015822,0.02090995573015823,50.644342039814,455.799078358326,0.002953923718372668,0.002953923718372669),
    @(3,6.228371,18.685113,0.1412687695991749,0.1412687695991749,3,243.3763986666667,730.1291960000001,0.625857000534647,0.6258570005346471,1515.838503539906,13642.54653185915,0.08841404841055971,0.08841404841055972),
    @(3,6.228371,18.685113,0.1412687695991749,0.1412687695991749,3,103.9426383333333,311.827915,0.2672947262403034,0.2672947262403035,647.3933142588216,5826.539828329395,0.03776039709631594,0.03776039709631595),
    @(3,6.228371,18.685113,0.1412687695991749,0.1412687695991749,3,33.41874933333333,100.256248,0.08593831749489127,0.08593831749489128,208.1443692040027,1873.299322836024,0.01214040037392653,0.01214040037392653)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $ws.Cells.Item($r, $cols[$j]).Value = $row[$j]
    }
}
